$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the rows that no longer have data (old rows 44-48) since the
# sector list shrank from 47 sectors to 42 sectors (plus header + Construction Materials).
$ws.Range("A44:B48").ClearContents()

# Update sector labels/values for rows 3-43 (row 2, "Construction Materials(8)", is unchanged)
$data = @(
    @(3, 'Gas Utilities(12)', 0.5186801690767942),
    @(4, 'Multi-Utilities(18)', 0.500909479822461),
    @(5, 'Marine(15)', 0.4852696261099415),
    @(6, 'Construction & Engineering(20)', 0.4342946658154804),
    @(7, 'Electric Utilities(28)', 0.4335275976640242),
    @(8, 'Building Products(23)', 0.423287851146404),
    @(9, 'Energy Equipment & Services(32)', 0.4014413576978428),
    @(10, 'Wireless Telecommunication Services(14)', 0.3941202375365044),
    @(11, 'Containers & Packaging(12)', 0.380250765115083),
    @(12, 'Capital Markets(75)', 0.3400709109274582),
    @(13, 'Auto Components(21)', 0.3381824331755282),
    @(14, 'Metals & Mining(89)', 0.3225174197229803),
    @(15, 'Life Sciences Tools & Services(19)', 0.3191484823244686),
    @(16, 'Chemicals(51)', 0.3147291236923744),
    @(17, 'Oil, Gas & Consumable Fuels(122)', 0.3014667954366632),
    @(18, 'Machinery(85)', 0.292889759709287),
    @(19, 'Insurance(75)', 0.2796246493324812),
    @(20, 'Diversified Consumer Services(17)', 0.2779035428496974),
    @(21, 'Diversified Telecommunication Services(20)', 0.2711477173287575),
    @(22, 'Pharmaceuticals(48)', 0.2676817266305845),
    @(23, 'Professional Services(35)', 0.2566271024674783),
    @(24, 'Banks(246)', 0.252534940170291),
    @(25, 'Communications Equipment(45)', 0.2349719760356817),
    @(26, 'Media(42)', 0.2298310009877373),
    @(27, 'Semiconductors & Semiconductor Equipment(68)', 0.2274804531565094),
    @(28, 'Household Durables(39)', 0.2268738898233296),
    @(29, 'Trading Companies & Distributors(25)', 0.2261234752372032),
    @(30, 'Real Estate Management & Development(22)', 0.2240819759548421),
    @(31, 'IT Services(52)', 0.2229508408562265),
    @(32, 'Health Care Providers & Services(46)', 0.2199865389515617),
    @(33, 'Road & Rail(22)', 0.199812233398616),
    @(34, 'Health Care Equipment & Supplies(83)', 0.1983760338025612),
    @(35, 'Software(66)', 0.1850055576448793),
    @(36, 'Thrifts & Mortgage Finance(47)', 0.184373386344324),
    @(37, 'Food Products(44)', 0.1839476891865025),
    @(38, 'Textiles, Apparel & Luxury Goods(29)', 0.178778919123491),
    @(39, 'Commercial Services & Supplies(52)', 0.1688693416197126),
    @(40, 'Hotels, Restaurants & Leisure(50)', 0.1511507946053877),
    @(41, 'Aerospace & Defense(37)', 0.1492547521842049),
    @(42, 'Specialty Retail(58)', 0.1408770635083056),
    @(43, 'Biotechnology(126)', 0.09277268615541284)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $label = $entry[1]
    $value = $entry[2]
    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $value
}
